# Applies the "Unnormalized Pmax Sample 2C 100kHz" edit:
#  - Column A (File Name) is replaced with the electrode location codes
#    that used to live in column C.
#  - Column B header changes from "Unnormalized P_max" to "P_max".
#  - Column A header changes from "File Name" to "Loc".
#  - Column C (Electrode Locations) is deleted entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Location codes, in row order (rows 2-64), taken from the old column C.
$locs = @(
    "A3","A5","A7","A11","A13","A15",
    "B9",
    "C1","C5","C13","C15",
    "D3","D7","D10","D12",
    "E1","E5","E7","E10","E13","E15",
    "F4","F10","F12",
    "G1","G5","G9","G12","G13",
    "H3","H7","H15",
    "I1","I3","I5","I8","I9","I11","I13","I15",
    "K1","K3","K5","K9","K11","K13","K15",
    "L7",
    "M1","M3","M6","M7","M9","M13","M15",
    "N12",
    "O1","O4","O7","O9","O11","O14","O15"
)

# Update column A with the location values (overwriting the filenames).
for ($i = 0; $i -lt $locs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $locs[$i]
}

# Delete column C (Electrode Locations) entirely, shifting nothing else -
# this removes the data and narrows the used range to A:B.
$ws.Columns.Item(3).Delete()

# Update header text.
$ws.Cells.Item(1, 1).Value = "Loc"
$ws.Cells.Item(1, 2).Value = "P_max"
